$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F21").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F23").Value = "ppe"
$ws.Range("F24").Value = "ppe"
$ws.Range("F28").Value = "application instructions"
$ws.Range("F29").Value = "env warning - water || off target movement"
$ws.Range("F31").Value = "env warning - species"
$ws.Range("F35").Value = "application instructions"
$ws.Range("F36").Value = "application instructions"
$ws.Range("F37").Value = "application instructions"
$ws.Range("F38").Value = "135_product_information"
$ws.Range("F39").Value = "135_product_information"
$ws.Range("F40").Value = "135_product_information"
$ws.Range("F41").Value = "135_product_information"
$ws.Range("F42").Value = "135_product_information"
$ws.Range("F43").Value = "135_product_information"
$ws.Range("F44").Value = "use restrictions"
$ws.Range("F46").Value = "use restrictions"
$ws.Range("F55").Value = "application instructions"
$ws.Range("F57").Value = "application instructions"
$ws.Range("F58").Value = "application instructions"
$ws.Range("F60").Value = "off target movement"
$ws.Range("F61").Value = "off target movement"
$ws.Range("F62").Value = "off target movement"
$ws.Range("F65").Value = "off target movement"
$ws.Range("F66").Value = "off target movement"
$ws.Range("F67").Value = "off target movement"
$ws.Range("F68").Value = "off target movement"
$ws.Range("F70").Value = "off target movement"
$ws.Range("F71").Value = "off target movement"
$ws.Range("F72").Value = "off target movement"
$ws.Range("F73").Value = "off target movement"
$ws.Range("F74").Value = "off target movement"
$ws.Range("F75").Value = "off target movement"
$ws.Range("F76").Value = "off target movement"
$ws.Range("F77").Value = "application instructions"
$ws.Range("F78").Value = "off target movement"
$ws.Range("F79").Value = "mixing"
$ws.Range("F81").Value = "mixing"
$ws.Range("F82").Value = "mixing"
$ws.Range("F83").Value = "mixing"
$ws.Range("F84").Value = "mixing"
$ws.Range("F85").Value = "mixing"
$ws.Range("F86").Value = "mixing"
$ws.Range("F87").Value = "mixing"
$ws.Range("F88").Value = "mixing"
$ws.Range("F89").Value = "mixing"
$ws.Range("F92").Value = "safety procedures"
$ws.Range("F93").Value = "safety procedures"
$ws.Range("F94").Value = "safety procedures"
$ws.Range("F95").Value = "safety procedures"
$ws.Range("F96").Value = "safety procedures"
$ws.Range("F97").Value = "safety procedures"
$ws.Range("F99").Value = "application instructions"
$ws.Range("F100").Value = "application instructions"
$ws.Range("F101").Value = "use restrictions"
$ws.Range("F102").Value = "use restrictions"
$ws.Range("F104").Value = "use restrictions"
$ws.Range("F105").Value = "use restrictions"
$ws.Range("F109").Value = "application instructions"
$ws.Range("F110").Value = "application instructions"
$ws.Range("F111").Value = "application instructions"
$ws.Range("F112").Value = "application instructions"
$ws.Range("F113").Value = "use restrictions"
$ws.Range("F141").Value = "mixing"
$ws.Range("F142").Value = "mixing"
$ws.Range("F143").Value = "mixing"
$ws.Range("F144").Value = "mixing"
$ws.Range("F145").Value = "mixing"
$ws.Range("F146").Value = "mixing"
$ws.Range("F149").Value = "mixing"
$ws.Range("F150").Value = "mixing"
$ws.Range("F151").Value = "mixing"
$ws.Range("F154").Value = "mixing"
$ws.Range("F156").Value = "154_pesticide_storage"
